$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A held full dates (2019-05-01 / 2019-10-01) but the header is "year";
# replace them with the plain text value "2019".
$ws.Range("A2").Value = "'2019"
$ws.Range("A3").Value = "'2019"
$ws.Range("A2:A3").ClearFormats()

# A trailing blank row (row 4) was added below the data.
$ws.Rows.Item(4).OutlineLevel = 0
